# Add a new "Sheet2" after the existing "Sheet1" and fill it with a summary
# of the Fiddler crab / Paraphronima model training data, then make it the
# active sheet (matching the author's "added summary of data" commit).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Fiddler crab training data summary ---
$ws2.Range("A1").Value = "Fiddler crab model training data"

$ws2.Range("A3").Value = "Train"
$ws2.Range("B3").Value = "Test"

$ws2.Range("A4").Value = "dampieri_20151218"
$ws2.Range("B4").Value = "dampieri_male_16"

$ws2.Range("A5").Value = "dampieri_20200218_male_left_1676"
$ws2.Range("B5").Value = "flammula_20190925_male_left"

$ws2.Range("A6").Value = "flammula_20180307"

$ws2.Range("A7").Value = "flammula_20200327_female_left_178_fullres_cropped"

# --- Paraphronima training data summary ---
$ws2.Range("A10").Value = "Paraphronima model training data"

$ws2.Range("A12").Value = "Train"
$ws2.Range("B12").Value = "Test"

# Column widths roughly matching the source workbook (46.4 / 18.38 chars,
# quantised to the nearest pixel Excel can actually store).
$ws2.Columns.Item(1).ColumnWidth = 45.5
$ws2.Columns.Item(2).ColumnWidth = 17.5

# Selections: Sheet1's cursor moves to C4, and the new Sheet2 (selection at
# B33) becomes the active/visible sheet.
[void]$ws1.Range("C4").Select()
[void]$ws2.Range("B33").Select()
